$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) / 1h-volume-change (E) figures for the cryptos list.
# Each cell must remain plain TEXT (it already holds things like "53.784.37"
# or "  -1.91%  " which are not valid numbers) so we briefly force a text
# number format before writing the value, then clear the format again so the
# cell style is left exactly as it was (General, unstyled).
$updates = @(
    @{ Cell = "D2"; Value = "53.784.37" }
    @{ Cell = "E2"; Value = "  -1.91%  " }
    @{ Cell = "D3"; Value = "2.253.74" }
    @{ Cell = "E3"; Value = "  -2.19%  " }
    @{ Cell = "E4"; Value = "  +0.25%  " }
    @{ Cell = "D5"; Value = "492.82" }
    @{ Cell = "E5"; Value = "  -1.17%  " }
    @{ Cell = "D6"; Value = "127.39" }
    @{ Cell = "E6"; Value = "  -1.25%  " }
    @{ Cell = "E7"; Value = "  +0.24%  " }
    @{ Cell = "D8"; Value = "0.524" }
    @{ Cell = "E8"; Value = "  -1.59%  " }
    @{ Cell = "D9"; Value = "0.0946" }
    @{ Cell = "E9"; Value = "  -0.99%  " }
    @{ Cell = "E10"; Value = "  +0.58%  " }
    @{ Cell = "E11"; Value = "  +2.63%  " }
    @{ Cell = "E12"; Value = "  +1.47%  " }
    @{ Cell = "D13"; Value = "2.660.33" }
    @{ Cell = "E13"; Value = "  -1.82%  " }
    @{ Cell = "D14"; Value = "22.51" }
    @{ Cell = "E14"; Value = "  +3.04%  " }
    @{ Cell = "D15"; Value = "53.777.83" }
    @{ Cell = "E15"; Value = "  -1.69%  " }
    @{ Cell = "E16"; Value = "  -0.73%  " }
    @{ Cell = "D17"; Value = "2.255.92" }
    @{ Cell = "E17"; Value = "  -1.02%  " }
    @{ Cell = "D18"; Value = "10.20" }
    @{ Cell = "E18"; Value = "  +0.61%  " }
    @{ Cell = "E19"; Value = "  +0.05%  " }
    @{ Cell = "D20"; Value = "301.40" }
    @{ Cell = "E20"; Value = "  -2.07%  " }
    @{ Cell = "D21"; Value = "6.27" }
    @{ Cell = "E21"; Value = "  -3.06%  " }
    @{ Cell = "D22"; Value = "0.999" }
    @{ Cell = "E22"; Value = "  +0.03%  " }
    @{ Cell = "D23"; Value = "60.81" }
    @{ Cell = "E23"; Value = "  -3.61%  " }
    @{ Cell = "E24"; Value = "  +0.12%  " }
    @{ Cell = "D25"; Value = "0.147" }
    @{ Cell = "E25"; Value = "  -2.93%  " }
    @{ Cell = "D26"; Value = "7.24" }
    @{ Cell = "E26"; Value = "  +1.18%  " }
    @{ Cell = "D27"; Value = "171.84" }
    @{ Cell = "E27"; Value = "  +0.69%  " }
    @{ Cell = "D28"; Value = "1.59" }
    @{ Cell = "E28"; Value = "  -1.13%  " }
    @{ Cell = "D29"; Value = "0.0₃0685" }
    @{ Cell = "E29"; Value = "  -2.11%  " }
    @{ Cell = "D30"; Value = "5.88" }
    @{ Cell = "E30"; Value = "  -0.80%  " }
    @{ Cell = "E31"; Value = "  -1.71%  " }
    @{ Cell = "E32"; Value = "  +0.00%  " }
    @{ Cell = "E33"; Value = "  -0.02%  " }
    @{ Cell = "D34"; Value = "0.998" }
    @{ Cell = "E34"; Value = "  -0.08%  " }
    @{ Cell = "D35"; Value = "0.929" }
    @{ Cell = "E35"; Value = "  +6.26%  " }
    @{ Cell = "E36"; Value = "  -1.17%  " }
    @{ Cell = "D37"; Value = "3.68" }
    @{ Cell = "E37"; Value = "  -0.47%  " }
    @{ Cell = "E38"; Value = "  -1.83%  " }
    @{ Cell = "E39"; Value = "  -2.26%  " }
    @{ Cell = "E40"; Value = "  -0.65%  " }
    @{ Cell = "D41"; Value = "124.36" }
    @{ Cell = "D42"; Value = "4.77" }
    @{ Cell = "E42"; Value = "  -1.77%  " }
    @{ Cell = "D43"; Value = "0.0487" }
    @{ Cell = "E43"; Value = "  -0.04%  " }
    @{ Cell = "D44"; Value = "0.0887" }
    @{ Cell = "E44"; Value = "  -0.83%  " }
    @{ Cell = "D45"; Value = "0.539" }
    @{ Cell = "E45"; Value = "  -2.29%  " }
    @{ Cell = "D46"; Value = "237.69" }
    @{ Cell = "E46"; Value = "  -2.98%  " }
    @{ Cell = "E47"; Value = "  -1.83%  " }
    @{ Cell = "D48"; Value = "0.0203" }
    @{ Cell = "E48"; Value = "  -0.21%  " }
    @{ Cell = "E49"; Value = "  +0.35%  " }
    @{ Cell = "D50"; Value = "16.03" }
    @{ Cell = "E50"; Value = "  -3.17%  " }
    @{ Cell = "D51"; Value = "4.62" }
    @{ Cell = "E51"; Value = "  -0.83%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
